$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("B2").Value = 0.40902777777777777
$ws.Range("C2").Value = 0.40972222222222227
$ws.Range("B3").Value = 0.41111111111111115
$ws.Range("C3").Value = 0.41319444444444442
$ws.Range("B4").Value = 0.41666666666666669
$ws.Range("C4").Value = 0.41736111111111113

$ws.Range("D8").Select()
